$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 424.625
$ws.Range("I2").Value = 279.4
$ws.Range("K2").Value = 279.4
$ws.Range("M2").Value = -166.4
# Row 4
$ws.Range("H4").Value = 583
$ws.Range("I4").Value = 583
$ws.Range("K4").Value = 583
$ws.Range("M4").Value = -469
# Row 15
$ws.Range("H15").Value = 731.2258
$ws.Range("I15").Value = 731.2258
$ws.Range("K15").Value = 2193.6774
$ws.Range("M15").Value = -2024.6774
# Row 33
$ws.Range("H33").Value = 135.33333
$ws.Range("J33").Value = 119.44444
$ws.Range("L33").Value = 119.44444
$ws.Range("N33").Value = -577.44444
# Row 58
$ws.Range("H58").Value = 1544.6364
$ws.Range("J58").Value = 3262.5
$ws.Range("L58").Value = 9787.5
$ws.Range("N58").Value = -10087.5
# Row 70
$ws.Range("H70").Value = 16012.857
$ws.Range("J70").Value = 21938
$ws.Range("L70").Value = 65814
$ws.Range("N70").Value = -66354
# Row 73
$ws.Range("H73").Value = 16012.857
$ws.Range("J73").Value = 21938
$ws.Range("L73").Value = 65814
$ws.Range("N73").Value = -67686
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
# Row 113
$ws.Range("H113").Value = 29000
$ws.Range("I113").Value = 29000
$ws.Range("K113").Value = 29000
$ws.Range("M113").Value = -25746
# Row 129
$ws.Range("H129").Value = 864.7463
$ws.Range("I129").Value = 632.6667
$ws.Range("K129").Value = 1898.0001
$ws.Range("M129").Value = 3101.9999
# Row 132
$ws.Range("H132").Value = 884.675
$ws.Range("I132").Value = 803.44116
$ws.Range("J132").Value = 1345
$ws.Range("K132").Value = 2410.32348
$ws.Range("L132").Value = 4035
$ws.Range("M132").Value = 119.67652
$ws.Range("N132").Value = -9095
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 141
$ws.Range("H141").Value = 4461.9287
$ws.Range("J141").Value = 6001.8
$ws.Range("L141").Value = 18005.4
$ws.Range("N141").Value = -28365.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2945.622
$ws.Range("I32").Value = 1777.1945
$ws.Range("J32").Value = 11358.3
$ws.Range("K32").Value = 1777.1945
$ws.Range("L32").Value = 11358.3
$ws.Range("M32").Value = -1490.1945
$ws.Range("N32").Value = -11932.3
# Row 45
$ws.Range("H45").Value = 1187.08
$ws.Range("I45").Value = 852.05884
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 852.05884
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -475.05884
$ws.Range("N45").Value = -2653
# Row 61
$ws.Range("H61").Value = 3712.0967
$ws.Range("I61").Value = 2810
$ws.Range("K61").Value = 2810
$ws.Range("M61").Value = -2598
# Row 131
$ws.Range("H131").Value = 44498.8
$ws.Range("J131").Value = 44498.8
$ws.Range("L131").Value = 44498.8
$ws.Range("N131").Value = -54578.8
# Row 132
$ws.Range("H132").Value = 1599.4286
$ws.Range("I132").Value = 1005.6875
$ws.Range("K132").Value = 3017.0625
$ws.Range("M132").Value = -487.0625
# Row 136
$ws.Range("H136").Value = 3712.0967
$ws.Range("I136").Value = 2810
$ws.Range("K136").Value = 8430
$ws.Range("M136").Value = -5880

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2360.0417
$ws.Range("I31").Value = 2039.1818
$ws.Range("K31").Value = 2039.1818
$ws.Range("M31").Value = -1744.1818
# Row 34
$ws.Range("H34").Value = 2360.0417
$ws.Range("I34").Value = 2039.1818
$ws.Range("K34").Value = 2039.1818
$ws.Range("M34").Value = -1837.1818
# Row 74
$ws.Range("H74").Value = 30624.5
$ws.Range("J74").Value = 30624.5
$ws.Range("L74").Value = 30624.5
$ws.Range("N74").Value = -32372.5
# Row 77
$ws.Range("H77").Value = 30624.5
$ws.Range("J77").Value = 30624.5
$ws.Range("L77").Value = 91873.5
$ws.Range("N77").Value = -100609.5
# Row 132
$ws.Range("H132").Value = 2394.5
$ws.Range("I132").Value = 1554.5883
$ws.Range("J132").Value = 3492.8462
$ws.Range("K132").Value = 4663.7649
$ws.Range("L132").Value = 10478.5386
$ws.Range("M132").Value = -2133.7649
$ws.Range("N132").Value = -15538.5386

$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 131
$ws.Range("H131").Value = 780.5
$ws.Range("I131").Value = 557.125
$ws.Range("J131").Value = 799.9239
$ws.Range("K131").Value = 1671.375
$ws.Range("L131").Value = 2399.7717
$ws.Range("M131").Value = 3368.625
$ws.Range("N131").Value = -12479.7717

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2173.0908
$ws.Range("I102").Value = 2173.0908
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2173.0908
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -551.0907999999999
$ws.Range("N102").ClearContents()
# Row 132
$ws.Range("H132").Value = 3499041.8
$ws.Range("I132").Value = 6412093
$ws.Range("K132").Value = 19236279
$ws.Range("M132").Value = -19233749

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3875.2144
$ws.Range("J7").Value = 9996.333000000001
$ws.Range("L7").Value = 9996.333000000001
$ws.Range("N7").Value = -10220.333
# Row 40
$ws.Range("H40").Value = 3572.2307
$ws.Range("I40").Value = 2231.3635
$ws.Range("J40").Value = 10947
$ws.Range("K40").Value = 2231.3635
$ws.Range("L40").Value = 10947
$ws.Range("M40").Value = -2095.3635
$ws.Range("N40").Value = -11219
# Row 122
$ws.Range("H122").Value = 4470.8096
$ws.Range("I122").Value = 3888.7896
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 11666.3688
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -9216.3688
$ws.Range("N122").Value = -34900
# Row 126
$ws.Range("H126").Value = 3875.2144
$ws.Range("J126").Value = 9996.333000000001
$ws.Range("L126").Value = 29988.999
$ws.Range("N126").Value = -34928.999
# Row 132
$ws.Range("H132").Value = 2764.5862
$ws.Range("I132").Value = 2094.0952
$ws.Range("K132").Value = 6282.285600000001
$ws.Range("M132").Value = -3752.285600000001

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2249.5
$ws.Range("I81").Value = 2166
$ws.Range("K81").Value = 4332
$ws.Range("M81").Value = -3271
# Row 84
$ws.Range("H84").Value = 2249.5
$ws.Range("I84").Value = 2166
$ws.Range("K84").Value = 21660
$ws.Range("M84").Value = -16356
# Row 107
$ws.Range("H107").Value = 897.6
$ws.Range("I107").Value = 662.1667
$ws.Range("J107").Value = 1250.75
$ws.Range("K107").Value = 1986.5001
$ws.Range("L107").Value = 3752.25
$ws.Range("M107").Value = -66.50009999999997
$ws.Range("N107").Value = -7592.25
# Row 132
$ws.Range("H132").Value = 1228.2307
$ws.Range("I132").Value = 845.619
$ws.Range("J132").Value = 2835.2
$ws.Range("K132").Value = 2536.857
$ws.Range("L132").Value = 8505.599999999999
$ws.Range("M132").Value = -6.856999999999971
$ws.Range("N132").Value = -13565.6
